$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new rows (bottom-to-top using original row numbers) to make room for new fields
$ws.Rows.Item(28).Insert()   # new row for t_surface_dataset.service_path
$ws.Rows.Item(11).Insert()   # new row for t_point_dataset.service_path
$ws.Rows.Item(4).Insert()    # new row for t_category.tenant_name

# Set final cell values for every populated row/column in the table
$ws.Range("A1").Value = '#'
$ws.Range("B1").Value = 'テーブル名'
$ws.Range("C1").Value = '列名'
$ws.Range("D1").Value = '型'
$ws.Range("E1").Value = '主キー'
$ws.Range("F1").Value = '制約'
$ws.Range("G1").Value = 'データシート列'

$ws.Range("B2").Value = 't_category'
$ws.Range("C2").Value = 'category_id'
$ws.Range("D2").Value = 'integer'
$ws.Range("E2").Value = '〇'
$ws.Range("F2").Value = 'not null unique'
$ws.Range("G2").Value = 'A'

$ws.Range("B3").Value = 't_category'
$ws.Range("C3").Value = 'category_name'
$ws.Range("D3").Value = 'varchar(50)'
$ws.Range("F3").Value = 'not null'
$ws.Range("G3").Value = 'B'

$ws.Range("B4").Value = 't_category'
$ws.Range("C4").Value = 'tenant_name'
$ws.Range("D4").Value = 'varchar(50)'
$ws.Range("F4").Value = 'not null'
$ws.Range("G4").Value = 'C'

$ws.Range("B5").Value = 't_category'
$ws.Range("C5").Value = 'category_color'
$ws.Range("D5").Value = 'varchar(32)'
$ws.Range("F5").Value = 'not null'
$ws.Range("G5").Value = 'D'

$ws.Range("B6").Value = 't_category'
$ws.Range("C6").Value = 'display_order'
$ws.Range("D6").Value = 'integer'
$ws.Range("F6").Value = 'not null'
$ws.Range("G6").Value = 'E'

$ws.Range("B7").Value = 't_category'
$ws.Range("C7").Value = 'enabled'
$ws.Range("D7").Value = 'boolean'
$ws.Range("F7").Value = 'not null'
$ws.Range("G7").Value = 'F'

$ws.Range("B9").Value = 't_point_dataset'
$ws.Range("C9").Value = 'point_dataset_id'
$ws.Range("D9").Value = 'integer'
$ws.Range("E9").Value = '〇'
$ws.Range("F9").Value = 'not null unique'
$ws.Range("G9").Value = 'A'

$ws.Range("B10").Value = 't_point_dataset'
$ws.Range("C10").Value = 'category_id'
$ws.Range("D10").Value = 'integer'
$ws.Range("F10").Value = 'references t_category(category_id) not null'
$ws.Range("G10").Value = 'B'

$ws.Range("B11").Value = 't_point_dataset'
$ws.Range("C11").Value = 'point_dataset_name'
$ws.Range("D11").Value = 'varchar(50)'
$ws.Range("F11").Value = 'not null'
$ws.Range("G11").Value = 'D'

$ws.Range("B12").Value = 't_point_dataset'
$ws.Range("C12").Value = 'service_path'
$ws.Range("D12").Value = 'varchar(550)'
$ws.Range("F12").Value = 'not null'
$ws.Range("G12").Value = 'E'

$ws.Range("B13").Value = 't_point_dataset'
$ws.Range("C13").Value = 'point_color_code'
$ws.Range("D13").Value = 'varchar(32)'
$ws.Range("F13").Value = 'not null'
$ws.Range("G13").Value = 'F'

$ws.Range("B14").Value = 't_point_dataset'
$ws.Range("C14").Value = 'entity_type'
$ws.Range("D14").Value = 'varchar(50)'
$ws.Range("F14").Value = 'not null'
$ws.Range("G14").Value = 'G'

$ws.Range("B15").Value = 't_point_dataset'
$ws.Range("C15").Value = 'coordinates_attr_name'
$ws.Range("D15").Value = 'varchar(50)'
$ws.Range("F15").Value = 'not null'
$ws.Range("G15").Value = 'I'

$ws.Range("B16").Value = 't_point_dataset'
$ws.Range("C16").Value = 'register_time_attr_name'
$ws.Range("D16").Value = 'varchar(50)'
$ws.Range("F16").Value = 'not null'
$ws.Range("G16").Value = 'J'

$ws.Range("B17").Value = 't_point_dataset'
$ws.Range("C17").Value = 'enabled'
$ws.Range("D17").Value = 'boolean'
$ws.Range("F17").Value = 'not null'
$ws.Range("G17").Value = 'H'

$ws.Range("B19").Value = 't_point_detail'
$ws.Range("C19").Value = 'point_detail_id'
$ws.Range("D19").Value = 'integer'
$ws.Range("E19").Value = '〇'
$ws.Range("F19").Value = 'not null unique'
$ws.Range("G19").Value = '#A'

$ws.Range("B20").Value = 't_point_detail'
$ws.Range("C20").Value = 'point_dataset_id'
$ws.Range("D20").Value = 'integer'
$ws.Range("F20").Value = 'references t_point_dataset(point_dataset_id) not null'
$ws.Range("G20").Value = 'A'

$ws.Range("B21").Value = 't_point_detail'
$ws.Range("C21").Value = 'item_attr_name'
$ws.Range("D21").Value = 'varchar(50)'
$ws.Range("F21").Value = 'not null'
$ws.Range("G21").Value = ':1'

$ws.Range("B22").Value = 't_point_detail'
$ws.Range("C22").Value = 'data_type'
$ws.Range("D22").Value = 'integer'
$ws.Range("F22").Value = 'not null'
$ws.Range("G22").Value = ':2'

$ws.Range("B23").Value = 't_point_detail'
$ws.Range("C23").Value = 'display_title'
$ws.Range("D23").Value = 'varchar(50)'
$ws.Range("F23").Value = 'not null'
$ws.Range("G23").Value = ':K*4'

$ws.Range("B24").Value = 't_point_detail'
$ws.Range("C24").Value = 'display_order'
$ws.Range("D24").Value = 'integer'
$ws.Range("F24").Value = 'not null'
$ws.Range("G24").Value = ':3'

$ws.Range("B25").Value = 't_point_detail'
$ws.Range("C25").Value = 'enabled'
$ws.Range("D25").Value = 'boolean'
$ws.Range("F25").Value = 'not null'
$ws.Range("G25").Value = 'H'

$ws.Range("B27").Value = 't_surface_dataset'
$ws.Range("C27").Value = 'surface_dataset_id'
$ws.Range("D27").Value = 'integer'
$ws.Range("E27").Value = '〇'
$ws.Range("F27").Value = 'not null unique'
$ws.Range("G27").Value = 'A'

$ws.Range("B28").Value = 't_surface_dataset'
$ws.Range("C28").Value = 'category_id'
$ws.Range("D28").Value = 'integer'
$ws.Range("F28").Value = 'references t_category(category_id) not null'
$ws.Range("G28").Value = 'B'

$ws.Range("B29").Value = 't_surface_dataset'
$ws.Range("C29").Value = 'surface_dataset_name'
$ws.Range("D29").Value = 'varchar(50)'
$ws.Range("F29").Value = 'not null'
$ws.Range("G29").Value = 'D'

$ws.Range("B30").Value = 't_surface_dataset'
$ws.Range("C30").Value = 'service_path'
$ws.Range("D30").Value = 'varchar(550)'
$ws.Range("F30").Value = 'not null'
$ws.Range("G30").Value = 'E'

$ws.Range("B31").Value = 't_surface_dataset'
$ws.Range("C31").Value = 'border_color_code'
$ws.Range("D31").Value = 'varchar(32)'
$ws.Range("F31").Value = 'not null'
$ws.Range("G31").Value = 'F'

$ws.Range("B32").Value = 't_surface_dataset'
$ws.Range("C32").Value = 'fill_color_code'
$ws.Range("D32").Value = 'varchar(32)'
$ws.Range("F32").Value = 'not null'
$ws.Range("G32").Value = 'G'

$ws.Range("B33").Value = 't_surface_dataset'
$ws.Range("C33").Value = 'entity_type'
$ws.Range("D33").Value = 'varchar(50)'
$ws.Range("F33").Value = 'not null'
$ws.Range("G33").Value = 'H'

$ws.Range("B34").Value = 't_surface_dataset'
$ws.Range("C34").Value = 'coordinates_attr_name'
$ws.Range("D34").Value = 'varchar(50)'
$ws.Range("F34").Value = 'not null'
$ws.Range("G34").Value = 'J'

$ws.Range("B35").Value = 't_surface_dataset'
$ws.Range("C35").Value = 'register_time_attr_name'
$ws.Range("D35").Value = 'varchar(50)'
$ws.Range("F35").Value = 'not null'
$ws.Range("G35").Value = 'K'

$ws.Range("B36").Value = 't_surface_dataset'
$ws.Range("C36").Value = 'enabled'
$ws.Range("D36").Value = 'boolean'
$ws.Range("F36").Value = 'not null'
$ws.Range("G36").Value = 'I'

$ws.Range("B38").Value = 't_surface_detail'
$ws.Range("C38").Value = 'surface_detail_id'
$ws.Range("D38").Value = 'integer'
$ws.Range("E38").Value = '〇'
$ws.Range("F38").Value = 'not null unique'
$ws.Range("G38").Value = '#A'

$ws.Range("B39").Value = 't_surface_detail'
$ws.Range("C39").Value = 'surface_dataset_id'
$ws.Range("D39").Value = 'integer'
$ws.Range("F39").Value = 'references t_surface_dataset(surface_dataset_id) not null'
$ws.Range("G39").Value = 'A'

$ws.Range("B40").Value = 't_surface_detail'
$ws.Range("C40").Value = 'item_attr_name'
$ws.Range("D40").Value = 'varchar(50)'
$ws.Range("F40").Value = 'not null'
$ws.Range("G40").Value = ':1'

$ws.Range("B41").Value = 't_surface_detail'
$ws.Range("C41").Value = 'display_title'
$ws.Range("D41").Value = 'varchar(50)'
$ws.Range("F41").Value = 'not null'
$ws.Range("G41").Value = ':L*3'

$ws.Range("B42").Value = 't_surface_detail'
$ws.Range("C42").Value = 'display_order'
$ws.Range("D42").Value = 'integer'
$ws.Range("F42").Value = 'not null'
$ws.Range("G42").Value = ':2'

$ws.Range("B43").Value = 't_surface_detail'
$ws.Range("C43").Value = 'enabled'
$ws.Range("D43").Value = 'boolean'
$ws.Range("F43").Value = 'not null'
$ws.Range("G43").Value = 'I'

# Match the final selection state from the diff
$ws.Range("G44").Select() | Out-Null
